# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# ranking sheet to the refreshed snapshot values.
#
# The source cells are stored as literal text (e.g. "283.05", "1.86%"),
# not numbers, so a plain `Range.Value = "..."` assignment would get
# auto-converted to a number/percentage by Excel. To preserve the exact
# text (and avoid leaving a stray NumberFormat on the cell), each cell is
# briefly switched to the "@" (Text) format while the value is written,
# then has its formatting cleared again so the cell's style matches the
# original (unstyled) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "283.05"
Set-TextValue "E2" "1.86%"

Set-TextValue "D3" "28.32"
Set-TextValue "E3" "3.93%"

Set-TextValue "D4" "5.017"
Set-TextValue "E4" "3.21%"

Set-TextValue "D5" "0.06508"
Set-TextValue "E5" "1.19%"

Set-TextValue "D6" "7.215"
Set-TextValue "E6" "2.87%"

Set-TextValue "D7" "3.351"

Set-TextValue "D8" "1.380"
Set-TextValue "E8" "15.91%"

Set-TextValue "D9" "0.9186"
Set-TextValue "E9" "3.39%"

Set-TextValue "D10" "0.1540"
Set-TextValue "E10" "-0.03%"

Set-TextValue "D11" "0.06609"
Set-TextValue "E11" "27.56%"

Set-TextValue "D12" "0.07553"
Set-TextValue "E12" "0.58%"

Set-TextValue "D13" "0.02809"
Set-TextValue "E13" "-2.58%"

Set-TextValue "D14" "0.08987"
Set-TextValue "E14" "0.11%"

Set-TextValue "D15" "0.001587"
Set-TextValue "E15" "1.32%"

Set-TextValue "D16" "0.0006340"
Set-TextValue "E16" "-0.44%"

Set-TextValue "D17" "0.006162"
Set-TextValue "E17" "0.43%"

Set-TextValue "D18" "3.442"
Set-TextValue "E18" "-1.03%"

Set-TextValue "D19" "2.238"
Set-TextValue "E19" "-1.44%"

Set-TextValue "D21" "0.1281"
Set-TextValue "E21" "-4.42%"

Set-TextValue "D22" "3.987"
Set-TextValue "E22" "1.74%"

Set-TextValue "E23" "1.77%"

Set-TextValue "D24" "0.04443"
Set-TextValue "E24" "0.74%"

Set-TextValue "D25" "0.001183"
Set-TextValue "E25" "0.54%"

Set-TextValue "E26" "14.07%"

Set-TextValue "D27" "0.0001199"
Set-TextValue "E27" "1.60%"

Set-TextValue "E28" "-1.57%"

Set-TextValue "D40" "0.04116"
Set-TextValue "E40" "-0.12%"

Set-TextValue "D41" "0.006691"
Set-TextValue "E41" "-1.84%"

Set-TextValue "D42" "0.1229"
Set-TextValue "E42" "4.71%"

Set-TextValue "D43" "0.002129"
Set-TextValue "E43" "11.41%"

Set-TextValue "D44" "0.01206"
Set-TextValue "E44" "3.75%"

Set-TextValue "E45" "6.16%"

Set-TextValue "D46" "1.966"
Set-TextValue "E46" "16.83%"
